$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 32   Number  27"
$ws.Range("C9").Value = "Report Covering the Week  6/30/2025  Through  7/6/2025"

# --- Weekly crime statistics table updates (rows 15-31) ---
$ws.Range("F15").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("H15").Value = -100
$ws.Range("N15").Value = -80
$ws.Range("C16").Value = 1
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = 42.857142857142
$ws.Range("I16").Value = 59
$ws.Range("J16").Value = 43
$ws.Range("K16").Value = 37.209302325581
$ws.Range("L16").Value = -6.349206349206
$ws.Range("M16").Value = -57.246376811594
$ws.Range("N16").Value = -93.099415204678
$ws.Range("C17").Value = 13
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = 85.714285714285
$ws.Range("F17").Value = 33
$ws.Range("G17").Value = 20
$ws.Range("H17").Value = 65
$ws.Range("I17").Value = 193
$ws.Range("J17").Value = 137
$ws.Range("K17").Value = 40.875912408759
$ws.Range("L17").Value = 17.682926829268
$ws.Range("M17").Value = 24.516129032258
$ws.Range("N17").Value = -59.368421052631
$ws.Range("D18").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("E18").Value = "'***.*"
$ws.Range("A14").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 37
$ws.Range("K18").Value = -2.631578947368
$ws.Range("L18").Value = 2.777777777777
$ws.Range("M18").Value = -75.657894736842
$ws.Range("N18").Value = -96.193415637860
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = 50
$ws.Range("F19").Value = 18
$ws.Range("G19").Value = 17
$ws.Range("H19").Value = 5.882352941176
$ws.Range("I19").Value = 155
$ws.Range("J19").Value = 126
$ws.Range("K19").Value = 23.015873015873
$ws.Range("L19").Value = -29.223744292237
$ws.Range("M19").Value = -23.645320197044
$ws.Range("N19").Value = -51.41065830721
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 1
$ws.Range("I14").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("E20").Value = 200
$ws.Range("K14").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("F20").Value = 13
$ws.Range("H20").Value = 30
$ws.Range("I20").Value = 47
$ws.Range("J20").Value = 41
$ws.Range("K20").Value = 14.634146341463
$ws.Range("L20").Value = -11.320754716981
$ws.Range("M20").Value = -42.682926829268
$ws.Range("N20").Value = -93.178519593613
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 14
$ws.Range("E21").Value = 78.571428571428
$ws.Range("F21").Value = 81
$ws.Range("G21").Value = 63
$ws.Range("H21").Value = 28.571428571428
$ws.Range("I21").Value = 500
$ws.Range("J21").Value = 402
$ws.Range("K21").Value = 24.378109452736
$ws.Range("L21").Value = -8.759124087591
$ws.Range("M21").Value = -32.34100135318
$ws.Range("N21").Value = -85.154394299287
$ws.Range("D22").Value = 1
$ws.Range("I14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("J22").Value = 7
$ws.Range("K22").Value = -85.714285714285
$ws.Range("M22").Value = -96.153846153846
$ws.Range("D24").Value = 9
$ws.Range("E24").Value = 88.888888888888
$ws.Range("F24").Value = 72
$ws.Range("H24").Value = 50
$ws.Range("I24").Value = 499
$ws.Range("J24").Value = 432
$ws.Range("K24").Value = 15.509259259259
$ws.Range("L24").Value = -10.412926391382
$ws.Range("M24").Value = 3.312629399585
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 3
$ws.Range("I14").Copy()
$ws.Range("D25").PasteSpecial(-4122)
$ws.Range("E25").Value = -33.333333333333
$ws.Range("K14").Copy()
$ws.Range("E25").PasteSpecial(-4122)
$ws.Range("F25").Value = 10
$ws.Range("G25").Value = 6
$ws.Range("H25").Value = 66.666666666666
$ws.Range("I25").Value = 131
$ws.Range("J25").Value = 78
$ws.Range("K25").Value = 67.948717948717
$ws.Range("L25").Value = 10.084033613445
$ws.Range("C26").Value = 9
$ws.Range("D26").Value = 17
$ws.Range("E26").Value = -47.058823529411
$ws.Range("F26").Value = 46
$ws.Range("G26").Value = 55
$ws.Range("H26").Value = -16.363636363636
$ws.Range("I26").Value = 253
$ws.Range("J26").Value = 283
$ws.Range("K26").Value = -10.600706713780
$ws.Range("L26").Value = -3.065134099616
$ws.Range("M26").Value = -39.473684210526
$ws.Range("D27").Value = 1
$ws.Range("I14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("F27").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = -100
$ws.Range("J27").Value = 16
$ws.Range("K27").Value = -50
$ws.Range("C28").Value = 1
$ws.Range("I14").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("D28").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").Value = "'***.*"
$ws.Range("A14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 11
$ws.Range("H28").Value = -63.636363636363
$ws.Range("I28").Value = 29
$ws.Range("K28").Value = -25.641025641025
$ws.Range("L28").Value = 3.571428571428
$ws.Range("F29").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("F29").PasteSpecial(-4122)
$ws.Range("N29").Value = -95.238095238095
$ws.Range("F30").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("F30").PasteSpecial(-4122)
$ws.Range("N30").Value = -95.161290322580
$ws.Range("D31").Value = 3
$ws.Range("I14").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E31").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E31").PasteSpecial(-4122)
$ws.Range("F31").Value = 2
$ws.Range("G31").Value = 3
$ws.Range("H31").Value = -33.333333333333
$ws.Range("I31").Value = 10
$ws.Range("J31").Value = 8
$ws.Range("K31").Value = 25
$ws.Range("L31").Value = 66.666666666666

